# Generate Report for Handback
# Applies the "handback" status update to the localization-status workbook:
#  - Overview sheet: status text "Ready for handoff" -> "Handed back: in sync with en-US"
#  - zh-cn / de-de sheets: status text updated, Latest Target File / Latest Handback File /
#    Latest Handback DateTime columns populated with handback results + hyperlinks
#  - Column widths widened to fit the new (longer) text

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"

$githubBase = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b8e257ce5b88600a14dc21e57bad446e3a3cd227/e2e/"
$zhHandbackDate = "2016-09-04 09:05:35"
$deHandbackDate = "2016-09-04 09:05:42"

$fileA = "0e620eba-fba9-4b6f-9f7b-0b6a9bccd57c"
$fileB = "99b528a1-3399-4f24-b7d6-b70c0dd02444"

$fileA_md = "$fileA.md"
$fileB_md = "$fileB.md"

# ---------------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("E2:F3").Value = $statusText

$wsOverview.Columns.Item(5).ColumnWidth = 29.166666666666668
$wsOverview.Columns.Item(6).ColumnWidth = 29.166666666666668

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("C2").Value = $statusText
$wsZh.Range("C3").Value = $statusText

$wsZh.Range("J2").Value = "$fileA.98e116e07aa3aab685e3b12626da22a9fc826387.zh-cn.xlf"
$wsZh.Range("K2").Value = $zhHandbackDate

$wsZh.Range("J3").Value = "$fileB.565a28dda4a268432ce14711d1922d7a58a177df.zh-cn.xlf"
$wsZh.Range("K3").Value = $zhHandbackDate

$wsZh.Hyperlinks.Add($wsZh.Range("I2"), "$githubBase$fileA_md", "", "", $fileA_md)
$wsZh.Hyperlinks.Add($wsZh.Range("I3"), "$githubBase$fileB_md", "", "", $fileB_md)

$wsZh.Columns.Item(3).ColumnWidth = 29.166666666666668
$wsZh.Columns.Item(9).ColumnWidth = 39.166666666666664
$wsZh.Columns.Item(10).ColumnWidth = 39.166666666666664

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("C2").Value = $statusText
$wsDe.Range("C3").Value = $statusText

$wsDe.Range("J2").Value = "$fileA.98e116e07aa3aab685e3b12626da22a9fc826387.de-de.xlf"
$wsDe.Range("K2").Value = $deHandbackDate

$wsDe.Range("J3").Value = "$fileB.565a28dda4a268432ce14711d1922d7a58a177df.de-de.xlf"
$wsDe.Range("K3").Value = $deHandbackDate

$wsDe.Hyperlinks.Add($wsDe.Range("I2"), "$githubBase$fileA_md", "", "", $fileA_md)
$wsDe.Hyperlinks.Add($wsDe.Range("I3"), "$githubBase$fileB_md", "", "", $fileB_md)

$wsDe.Columns.Item(3).ColumnWidth = 29.166666666666668
$wsDe.Columns.Item(9).ColumnWidth = 39.166666666666664
$wsDe.Columns.Item(10).ColumnWidth = 39.166666666666664

Write-Output "Handback report applied"
